# Section Properties workbook update
# - Add "Top Stringer1" (row 6) z'-position data (E6), which previously held
#   a "TBD" placeholder string and produced #VALUE! errors downstream.
# - Recompute "Bottom Stringer1"/"Bottom Stringer2" (rows 11-12) y'/z'
#   positions (D/E) from updated geometry formulas.
# - Extend the section totals (row 16) SUM ranges to include the now
#   populated row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Top Stringer1 (row 6): replace "TBD" placeholder with real data ---
$ws.Range("E6").Value = -0.0625
$ws.Range("E6").Style = "Normal"

# --- Bottom Stringer1 (row 11): new y'/z' position formulas ---
$ws.Range("D11").Formula = "=-1.6046754518"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "=-0.092597515"
$ws.Range("E11").Style = "Normal"

# --- Bottom Stringer2 (row 12): new y'/z' position formulas ---
$ws.Range("D12").Formula = "=-1.7254548192"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "=-0.1248053464"
$ws.Range("E12").Style = "Normal"

# --- TOTALS row (16): widen SUM ranges from C2:C5 to C2:C6, etc. ---
$ws.Range("C16").Formula  = "=SUM(C2:C6,C11:C12)"
$ws.Range("F16").Formula  = "=SUM(F2:F6,F11:F12)"
$ws.Range("G16").Formula  = "=SUM(G2:G6,G11:G12)"
$ws.Range("H16").Formula  = "=SUM(H2:H6,H11:H12)"
$ws.Range("I16").Formula  = "=SUM(I2:I6,I11:I12)"
$ws.Range("J16").Formula  = "=SUM(J2:J6,J11:J12)"
$ws.Range("M16").Formula  = "=SUM(M2:M6,M11:M12)"
$ws.Range("N16").Formula  = "=SUM(N2:N6,N11:N12)"
$ws.Range("O16").Formula  = "=SUM(O2:O6,O11:O12)"
$ws.Range("R16").Formula  = "=SUM(R2:R6,R11:R12)"
$ws.Range("S16").Formula  = "=SUM(S2:S6,S11:S12)"
$ws.Range("T16").Formula  = "=SUM(T2:T6,T11:T12)"
$ws.Range("W16").Formula  = "=SUM(W2:W6,W11:W12)"
$ws.Range("X16").Formula  = "=SUM(X2:X6,X11:X12)"
$ws.Range("Y16").Formula  = "=SUM(Y2:Y6,Y11:Y12)"
$ws.Range("AZ16").Formula = "=SUM(AZ2:AZ6,AZ11:AZ12)"
$ws.Range("BA16").Formula = "=SUM(BA2:BA6,BA11:BA12)"
$ws.Range("BB16").Formula = "=SUM(BB2:BB6,BB11:BB12)"

# Row 6 formulas depended on E6 while it still held the "TBD" placeholder
# and were cached as #VALUE! errors. Re-apply their (unchanged) formula
# text so the engine re-evaluates them against the new numeric E6.
$ws.Range("G6").Formula  = "=C6*E6"
$ws.Range("L6").Formula  = "=E6-`$AI`$3"
$ws.Range("N6").Formula  = "=C6*L6^2"
$ws.Range("O6").Formula  = "=C6*K6*L6"
$ws.Range("T6").Formula  = "=R6*E6"
$ws.Range("V6").Formula  = "=E6-`$AI`$5"
$ws.Range("X6").Formula  = "=C6*V6^2"
$ws.Range("Y6").Formula  = "=C6*U6*V6"
$ws.Range("AZ6").Formula = "=(Q6)*(H6+X6)"
$ws.Range("BB6").Formula = "=(Q6)*(J6+Y6)"
